{"js": "// Update the date line and every arithmetic-problem cell in the single\n// table, preserving existing run formatting (Arial/TimeNewRoman, sz 30).\n\nconst NEW_DATE = \"2025-07-18 Friday\";\nconst NEW_VALUES = [\n  [\"36+46=\", \"40-14=\", \"40-8=\", \"36+59=\", \"53-44=\"],\n  [\"44+27=\", \"81-8=\", \"13-9=\", \"29+64=\", \"7+7=\"],\n  [\"56-38=\", \"16+35=\", \"84-49=\", \"62-19=\", \"73-24=\"],\n  [\"37+39=\", \"27+19=\", \"19+22=\", \"70-55=\", \"90-39=\"],\n  [\"34-5=\", \"83+8=\", \"74-27=\", \"86+5=\", \"37-8=\"],\n  [\"35+17=\", \"42-18=\", \"53-19=\", \"28+7=\", \"70-37=\"],\n  [\"57-39=\", \"25+17=\", \"74+8=\", \"72-3=\", \"37+54=\"],\n  [\"67+24=\", \"61-44=\", \"52-19=\", \"17+56=\", \"70-17=\"],\n  [\"27+44=\", \"94-27=\", \"92-58=\", \"30-25=\", \"85-16=\"],\n  [\"44+37=\", \"83-76=\", \"17+54=\", \"14+79=\", \"76-7=\"],\n  [\"17+28=\", \"18+79=\", \"16+6=\", \"34-9=\", \"21-7=\"],\n  [\"82-38=\", \"67+14=\", \"66-7=\", \"49+43=\", \"45-38=\"],\n  [\"9+16=\", \"43+48=\", \"59+28=\", \"97-59=\", \"49+37=\"],\n  [\"24-8=\", \"24-9=\", \"51-7=\", \"45+6=\", \"82-25=\"],\n  [\"37+45=\", \"64-29=\", \"19+19=\", \"32-16=\", \"69+25=\"],\n  [\"67+26=\", \"22-3=\", \"81-12=\", \"64-45=\", \"23+49=\"],\n  [\"51-14=\", \"62-25=\", \"84-26=\", \"45-8=\", \"47+24=\"],\n  [\"6+56=\", \"29+64=\", \"70-59=\", \"7+86=\", \"8+4=\"],\n  [\"4+37=\", \"69+27=\", \"41-18=\", \"85-8=\", \"36+7=\"],\n  [\"76-57=\", \"13+49=\", \"6+45=\", \"24+27=\", \"68+6=\"],\n];\n\n// 1. Update the title/date paragraph (first paragraph in the body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(NEW_DATE, Word.InsertLocation.replace);\n\n// 2. Update every cell of the (single) table in document order.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\ntable.values = NEW_VALUES;\n\nawait context.sync();\n", "ps1": "# Update the date line and every arithmetic-problem cell in the single\n# table, preserving existing run formatting (Arial/TimeNewRoman, sz 30).\n\n$d = $word.ActiveDocument\n\n# 1. Update the title/date paragraph (first paragraph in the body).\n$d.Paragraphs.Item(1).Range.Text = \"2025-07-18 Friday\"\n\n# 2. Update every cell of the (single) table in document order.\n$newValues = @(\n    @(\"36+46=\", \"40-14=\", \"40-8=\", \"36+59=\", \"53-44=\"),\n    @(\"44+27=\", \"81-8=\", \"13-9=\", \"29+64=\", \"7+7=\"),\n    @(\"56-38=\", \"16+35=\", \"84-49=\", \"62-19=\", \"73-24=\"),\n    @(\"37+39=\", \"27+19=\", \"19+22=\", \"70-55=\", \"90-39=\"),\n    @(\"34-5=\", \"83+8=\", \"74-27=\", \"86+5=\", \"37-8=\"),\n    @(\"35+17=\", \"42-18=\", \"53-19=\", \"28+7=\", \"70-37=\"),\n    @(\"57-39=\", \"25+17=\", \"74+8=\", \"72-3=\", \"37+54=\"),\n    @(\"67+24=\", \"61-44=\", \"52-19=\", \"17+56=\", \"70-17=\"),\n    @(\"27+44=\", \"94-27=\", \"92-58=\", \"30-25=\", \"85-16=\"),\n    @(\"44+37=\", \"83-76=\", \"17+54=\", \"14+79=\", \"76-7=\"),\n    @(\"17+28=\", \"18+79=\", \"16+6=\", \"34-9=\", \"21-7=\"),\n    @(\"82-38=\", \"67+14=\", \"66-7=\", \"49+43=\", \"45-38=\"),\n    @(\"9+16=\", \"43+48=\", \"59+28=\", \"97-59=\", \"49+37=\"),\n    @(\"24-8=\", \"24-9=\", \"51-7=\", \"45+6=\", \"82-25=\"),\n    @(\"37+45=\", \"64-29=\", \"19+19=\", \"32-16=\", \"69+25=\"),\n    @(\"67+26=\", \"22-3=\", \"81-12=\", \"64-45=\", \"23+49=\"),\n    @(\"51-14=\", \"62-25=\", \"84-26=\", \"45-8=\", \"47+24=\"),\n    @(\"6+56=\", \"29+64=\", \"70-59=\", \"7+86=\", \"8+4=\"),\n    @(\"4+37=\", \"69+27=\", \"41-18=\", \"85-8=\", \"36+7=\"),\n    @(\"76-57=\", \"13+49=\", \"6+45=\", \"24+27=\", \"68+6=\")\n)\n\n$table = $d.Tables.Item(1)\nfor ($r = 1; $r -le $newValues.Count; $r++) {\n    $row = $newValues[$r - 1]\n    for ($c = 1; $c -le $row.Count; $c++) {\n        $table.Cell($r, $c).Range.Text = $row[$c - 1]\n    }\n}\n"}
